$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where coin name (B) and link (C) were reordered/swapped
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("B39").Value = "NEARProtocol"
$ws.Range("C39").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("B44").Value = "Algorand"
$ws.Range("C44").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("B45").Value = "BitcoinSV"
$ws.Range("C45").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"

# Price (D) and Volume(1h) (E) updates
# D column holds text-formatted numbers (e.g. "318.02", "12.00"); force Text
# number format first so Excel does not auto-convert them to numeric values.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.368.79"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.290.70"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.02"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.37"
$ws.Range("E6").Value = "  -4.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.628"
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  +0.17%  "
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.27"
$ws.Range("E10").Value = "  -2.43%  "
$ws.Range("E11").Value = "  -1.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.27"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.107"
$ws.Range("E13").Value = "  +0.24%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.963"
$ws.Range("E14").Value = "  -1.85%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.21"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.639.00"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.295.03"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.331.34"
$ws.Range("E18").Value = "  +0.45%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.42"
$ws.Range("E19").Value = "  -3.43%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.92"
$ws.Range("E21").Value = "  +30.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.92"
$ws.Range("E22").Value = "  -0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.56"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.73"
$ws.Range("E24").Value = "  +2.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.22"
$ws.Range("E25").Value = "  -5.21%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.84"
$ws.Range("E27").Value = "  -1.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.33"
$ws.Range("E28").Value = "  +2.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "22.59"
$ws.Range("E29").Value = "  -1.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.59"
$ws.Range("E30").Value = "  +2.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "166.19"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.13"
$ws.Range("E32").Value = "  +3.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0873"
$ws.Range("E33").Value = "  -2.45%  "
$ws.Range("E34").Value = "  +1.93%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.116"
$ws.Range("E35").Value = "  -4.31%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.55"
$ws.Range("E36").Value = "  -12.81%  "
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0357"
$ws.Range("E38").Value = "  +0.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.66"
$ws.Range("E39").Value = "  +0.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.74"
$ws.Range("E40").Value = "  -7.77%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.53"
$ws.Range("E41").Value = "  +2.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "68.69"
$ws.Range("E42").Value = "  -3.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  -0.10%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.224"
$ws.Range("E44").Value = "  -1.94%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "93.21"
$ws.Range("E45").Value = "  -7.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "115.31"
$ws.Range("E46").Value = "  +1.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.00"
$ws.Range("E47").Value = "  -2.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "78.89"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.97"
$ws.Range("E49").Value = "  -2.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.23"
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.604.83"
$ws.Range("E51").Value = "  +3.48%  "
